# IRYO-vaccination_data.xlsx — add the 2021-05-20 daily entry
#
# The published sheet is a "newest date on top" running log (row 4 holds the
# running grand totals, row 5 is always the most-recent day). This update
# adds one new day row (2021-05-20, a Thursday) above the previous top data
# row, pushing every existing day row down by one, and refreshes the
# cumulative totals in row 4 plus the "as of" caption to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row above the current first data row (row 5), shifting
#    all existing day rows (previously 5-28) down to (6-29).
$ws.Rows("5:5").Insert()

# Row-insert in this host leaves the new row with a generic default format;
# restore it to look exactly like the data rows around it (same date /
# centered-day-label / number styles used throughout the table) by copying
# the formatting back from the row immediately below (the old row 5, now
# row 6).
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Populate the new top row with the 2021-05-20 figures.
$ws.Range("A5").Value2 = 44336          # 2021-05-20 (serial date)
$ws.Range("B5").Value2 = "(木)"         # Thursday label (matches existing shared text)
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value2 = 81422
$ws.Range("E5").Value2 = 109111

# 3) Refresh the cumulative-total row (row 4). C4 is already =SUM(D4:E4),
#    so it recalculates once D4/E4 are updated.
$ws.Range("D4").Value2 = 3865493
$ws.Range("E4").Value2 = 2323873

# 4) Update the "as of" caption from 5/19 to 5/20.
$ws.Range("E2").Value2 = "（5月20日時点）"
